# Apply updated leve profit calculations across ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 132
$wsALC.Range("H132").Value = 974
$wsALC.Range("I132").Value = 824.6667
$wsALC.Range("K132").Value = 2474.0001
$wsALC.Range("M132").Value = 55.9998999999998

# ALC row 137
$wsALC.Range("H137").Value = 3410.3333
$wsALC.Range("I137").Value = 2719.5
$wsALC.Range("K137").Value = 8158.5
$wsALC.Range("M137").Value = -5608.5

# ALC row 138
$wsALC.Range("H138").Value = 3390.7827
$wsALC.Range("J138").Value = 4318.6523
$wsALC.Range("L138").Value = 12955.9569
$wsALC.Range("N138").Value = -23235.9569

# ARM row 2
$wsARM.Range("H2").Value = 34666.332
$wsARM.Range("I2").Value = 1999.5
$wsARM.Range("K2").Value = 1999.5
$wsARM.Range("M2").Value = -1886.5

# ARM row 43
$wsARM.Range("H43").Value = 30797.889
$wsARM.Range("I43").Value = 29171
$wsARM.Range("K43").Value = 29171
$wsARM.Range("M43").Value = -28858

# ARM row 45
$wsARM.Range("H45").Value = 52634330
$wsARM.Range("I45").Value = 71429910
$wsARM.Range("J45").Value = 6704.4
$wsARM.Range("K45").Value = 71429910
$wsARM.Range("L45").Value = 6704.4
$wsARM.Range("M45").Value = -71429533
$wsARM.Range("N45").Value = -7458.4

# ARM row 61
$wsARM.Range("H61").Value = 5699.4443
$wsARM.Range("I61").Value = 4776.975
$wsARM.Range("K61").Value = 4776.975
$wsARM.Range("M61").Value = -4564.975

# ARM row 62
$wsARM.Range("H62").Value = 40000
$wsARM.Range("J62").Value = 40000
$wsARM.Range("L62").Value = 40000
$wsARM.Range("N62").Value = -41248

# ARM row 65
$wsARM.Range("H65").Value = 40000
$wsARM.Range("J65").Value = 40000
$wsARM.Range("L65").Value = 120000
$wsARM.Range("N65").Value = -126240

# ARM row 74
$wsARM.Range("H74").Value = 12347275
$wsARM.Range("I74").Value = 14494361
$wsARM.Range("K74").Value = 14494361
$wsARM.Range("M74").Value = -14493487

# ARM row 77
$wsARM.Range("H77").Value = 12347275
$wsARM.Range("I77").Value = 14494361
$wsARM.Range("K77").Value = 72471805
$wsARM.Range("M77").Value = -72467437

# ARM row 102
$wsARM.Range("H102").Value = 2439.647
$wsARM.Range("I102").Value = 2404.625
$wsARM.Range("K102").Value = 2404.625
$wsARM.Range("M102").Value = -782.625

# ARM row 116
$wsARM.Range("H116").Value = 34666.332
$wsARM.Range("I116").Value = 1999.5
$wsARM.Range("K116").Value = 1999.5
$wsARM.Range("M116").Value = 294.5

# ARM row 122
$wsARM.Range("H122").Value = 2142.5715
$wsARM.Range("I122").Value = 1312.5714
$wsARM.Range("K122").Value = 3937.7142
$wsARM.Range("M122").Value = -1487.7142

# ARM row 132
$wsARM.Range("H132").Value = 2103.697
$wsARM.Range("I132").Value = 1743.661
$wsARM.Range("K132").Value = 5230.983
$wsARM.Range("M132").Value = -2700.983

# ARM row 136
$wsARM.Range("H136").Value = 5699.4443
$wsARM.Range("I136").Value = 4776.975
$wsARM.Range("K136").Value = 14330.925
$wsARM.Range("M136").Value = -11780.925

# BSM row 3
$wsBSM.Range("H3").Value = 34666.332
$wsBSM.Range("I3").Value = 1999.5
$wsBSM.Range("K3").Value = 1999.5
$wsBSM.Range("M3").Value = -1885.5

# BSM row 20
$wsBSM.Range("H20").Value = 2656.842
$wsBSM.Range("I20").Value = 1742.7142
$wsBSM.Range("J20").Value = 5216.4
$wsBSM.Range("K20").Value = 1742.7142
$wsBSM.Range("L20").Value = 5216.4
$wsBSM.Range("M20").Value = -1495.7142
$wsBSM.Range("N20").Value = -5710.4

# BSM row 99
$wsBSM.Range("H99").Value = 1265.8334
$wsBSM.Range("I99").Value = 1127
$wsBSM.Range("K99").Value = 1127
$wsBSM.Range("M99").Value = 371

# BSM row 107
$wsBSM.Range("H107").Value = 2899.625
$wsBSM.Range("I107").Value = 2599.5715
$wsBSM.Range("K107").Value = 2599.5715
$wsBSM.Range("M107").Value = -679.5715

# BSM row 140
$wsBSM.Range("H140").Value = 69999
$wsBSM.Range("J140").Value = 69999
$wsBSM.Range("L140").Value = 69999
$wsBSM.Range("N140").Value = -80359

# CRP row 7
$wsCRP.Range("H7").Value = 212
$wsCRP.Range("I7").Value = 44.5
$wsCRP.Range("K7").Value = 44.5
$wsCRP.Range("M7").Value = 68.5

# CRP row 16
$wsCRP.Range("H16").Value = 2476.8
$wsCRP.Range("I16").Value = 1647.25
$wsCRP.Range("K16").Value = 1647.25
$wsCRP.Range("M16").Value = -1360.25

# CRP row 22
$wsCRP.Range("H22").Value = 1659.0667
$wsCRP.Range("I22").Value = 357.14285
$wsCRP.Range("K22").Value = 357.14285
$wsCRP.Range("M22").Value = -7.14285000000001

# CRP row 113
$wsCRP.Range("H113").Value = 2476.8
$wsCRP.Range("I113").Value = 1647.25
$wsCRP.Range("K113").Value = 1647.25
$wsCRP.Range("M113").Value = 522.75

# CRP row 132
$wsCRP.Range("H132").Value = 21664.666
$wsCRP.Range("I132").Value = 21664.666
$wsCRP.Range("K132").Value = 64993.99800000001
$wsCRP.Range("M132").Value = -62463.99800000001

# GSM row 19
$wsGSM.Range("H19").Value = 0
$wsGSM.Range("I19").Value = 0
$wsGSM.Range("K19").Value = 0
$wsGSM.Range("M19").Value = $null

# GSM row 21
$wsGSM.Range("H21").Value = 30003
$wsGSM.Range("I21").Value = 30003
$wsGSM.Range("J21").Value = 0
$wsGSM.Range("K21").Value = 30003
$wsGSM.Range("L21").Value = 0
$wsGSM.Range("M21").Value = -29830
$wsGSM.Range("N21").Value = $null

# GSM row 22
$wsGSM.Range("H22").Value = 9500
$wsGSM.Range("I22").Value = 0
$wsGSM.Range("J22").Value = 9500
$wsGSM.Range("K22").Value = 0
$wsGSM.Range("L22").Value = 9500
$wsGSM.Range("M22").Value = $null
$wsGSM.Range("N22").Value = -10558

# GSM row 30
$wsGSM.Range("H30").Value = 30003
$wsGSM.Range("I30").Value = 30003
$wsGSM.Range("J30").Value = 0
$wsGSM.Range("K30").Value = 30003
$wsGSM.Range("L30").Value = 0
$wsGSM.Range("M30").Value = -29898
$wsGSM.Range("N30").Value = $null

# GSM row 53
$wsGSM.Range("H53").Value = 0
$wsGSM.Range("J53").Value = 0
$wsGSM.Range("L53").Value = 0
$wsGSM.Range("N53").Value = $null

# GSM row 70
$wsGSM.Range("H70").Value = 7471
$wsGSM.Range("I70").Value = 5434.5
$wsGSM.Range("K70").Value = 5434.5
$wsGSM.Range("M70").Value = -5164.5

# GSM row 73
$wsGSM.Range("H73").Value = 7471
$wsGSM.Range("I73").Value = 5434.5
$wsGSM.Range("K73").Value = 5434.5
$wsGSM.Range("M73").Value = -4498.5

# GSM row 113
$wsGSM.Range("H113").Value = 3920.2632
$wsGSM.Range("I113").Value = 3092.4614
$wsGSM.Range("J113").Value = 5713.8335
$wsGSM.Range("K113").Value = 3092.4614
$wsGSM.Range("L113").Value = 5713.8335
$wsGSM.Range("M113").Value = -922.4614000000001
$wsGSM.Range("N113").Value = -10053.8335

# GSM row 122
$wsGSM.Range("H122").Value = 4503.143
$wsGSM.Range("I122").Value = 3976.7368
$wsGSM.Range("K122").Value = 11930.2104
$wsGSM.Range("M122").Value = -9480.2104

# LTW row 62
$wsLTW.Range("H62").Value = 255000
$wsLTW.Range("J62").Value = 255000
$wsLTW.Range("L62").Value = 255000
$wsLTW.Range("N62").Value = -256248

# LTW row 65
$wsLTW.Range("H65").Value = 255000
$wsLTW.Range("J65").Value = 255000
$wsLTW.Range("L65").Value = 765000
$wsLTW.Range("N65").Value = -771240

# LTW row 104
$wsLTW.Range("H104").Value = 30092
$wsLTW.Range("J104").Value = 30092
$wsLTW.Range("L104").Value = 30092
$wsLTW.Range("N104").Value = -37080

# LTW row 132
$wsLTW.Range("H132").Value = 9376.25
$wsLTW.Range("I132").Value = 3625
$wsLTW.Range("J132").Value = 11293.333
$wsLTW.Range("K132").Value = 10875
$wsLTW.Range("L132").Value = 33879.999
$wsLTW.Range("M132").Value = -8345
$wsLTW.Range("N132").Value = -38939.999

# WVR row 19
$wsWVR.Range("H19").Value = 550
$wsWVR.Range("I19").Value = 100
$wsWVR.Range("J19").Value = 1000
$wsWVR.Range("K19").Value = 100
$wsWVR.Range("L19").Value = 1000
$wsWVR.Range("M19").Value = 74
$wsWVR.Range("N19").Value = -1348

# WVR row 104
$wsWVR.Range("H104").Value = 7500
$wsWVR.Range("J104").Value = 7500
$wsWVR.Range("L104").Value = 7500
$wsWVR.Range("N104").Value = -14488

# WVR row 108
$wsWVR.Range("H108").Value = 30000
$wsWVR.Range("J108").Value = 30000
$wsWVR.Range("L108").Value = 30000
$wsWVR.Range("N108").Value = -37680

# WVR row 121
$wsWVR.Range("H121").Value = 68421
$wsWVR.Range("J121").Value = 68421
$wsWVR.Range("L121").Value = 68421
$wsWVR.Range("N121").Value = -71915

# WVR row 136
$wsWVR.Range("H136").Value = 4125.3335
$wsWVR.Range("I136").Value = 2301.647
$wsWVR.Range("J136").Value = 11876
$wsWVR.Range("K136").Value = 6904.941
$wsWVR.Range("L136").Value = 35628
$wsWVR.Range("M136").Value = -4354.941
$wsWVR.Range("N136").Value = -40728
